$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 160, shifting existing rows 160-176 down to 161-177
$ws.Rows("160:160").Insert()

# Populate the newly inserted row 160 with the new data record
$ws.Cells.Item(160, 1).Value = 9
$ws.Cells.Item(160, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(160, 3).Value = "Metropolitana"
$ws.Cells.Item(160, 4).Value = 45258
$ws.Cells.Item(160, 5).Value = 13
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100101
$ws.Cells.Item(160, 8).Value = "Berries"
$ws.Cells.Item(160, 9).Value = 100101004
$ws.Cells.Item(160, 10).Value = "Frambuesa"
$ws.Cells.Item(160, 11).Value = "Sin especificar"
$ws.Cells.Item(160, 12).Value = "Primera"
$ws.Cells.Item(160, 13).Value = 120
$ws.Cells.Item(160, 14).Value = 11000
$ws.Cells.Item(160, 15).Value = 11000
$ws.Cells.Item(160, 16).Value = 11000
$ws.Cells.Item(160, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(160, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(160, 19).Value = 5500
$ws.Cells.Item(160, 20).Value = 2

# Ensure date styling/format matches the other rows in column D
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(161, 4).NumberFormat
